# EIA Table 2.2.B monthly update: October 2016 -> November 2016
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Update title / header text (shared strings content)
# ---------------------------------------------------------------
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Thousand Barrels)"
$ws.Range("A57").Value = "Rolling 12 Months Ending in November"

# ---------------------------------------------------------------
# 2. Insert a new row 53 for the "November" Year-to-Date entry.
#    This shifts the old rows 53-60 down to 54-61 and automatically
#    fixes up the sheet dimension and merged-cell ranges.
# ---------------------------------------------------------------
$ws.Rows(53).Insert()

# Copy the formatting from the October row (row 52) onto the new
# November row so it keeps the same styles (no new/duplicate styles).
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 184
$ws.Range("C53").Value = 0.35
$ws.Range("D53").Value = 78
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 97

# ---------------------------------------------------------------
# 3. Update the "Year to Date" block (now rows 55-57).
# ---------------------------------------------------------------
# Year 2014
$ws.Range("B55").Value = 2899
$ws.Range("C55").Value = 63
$ws.Range("D55").Value = 1089
$ws.Range("E55").Value = 212
$ws.Range("F55").Value = 1536

# Year 2015
$ws.Range("A56").Value = 2015
$ws.Range("B56").Value = 2932
$ws.Range("C56").Value = 61
$ws.Range("D56").Value = 1065
$ws.Range("E56").Value = 277
$ws.Range("F56").Value = 1529

# Year 2016
$ws.Range("A57").Value = 2016
$ws.Range("B57").Value = 2307
$ws.Range("C57").Value = 17
$ws.Range("D57").Value = 919
$ws.Range("E57").Value = 113
$ws.Range("F57").Value = 1259

# ---------------------------------------------------------------
# 4. Update the "Rolling 12 Months" block (now rows 59-60).
# ---------------------------------------------------------------
# Ending 2015
$ws.Range("B59").Value = 3132
$ws.Range("C59").Value = 62
$ws.Range("D59").Value = 1146
$ws.Range("E59").Value = 281
$ws.Range("F59").Value = 1643

# Ending 2016
$ws.Range("B60").Value = 2517
$ws.Range("C60").Value = 18
$ws.Range("D60").Value = 1009
$ws.Range("F60").Value = 1372
